$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert two new trigger rows ("L1_DoubleMu0" and "L1_DoubleMu0_SQ")
# right before the existing row 3 ("L1_SingleMuOpen"). This pushes
# the existing rows 3 and 4 down to rows 5 and 6.
# ------------------------------------------------------------------
$ws.Range("A3:A4").EntireRow.Insert()

# Copy the number format / style of the existing data columns (C:K)
# onto the two freshly inserted rows so they share the same style
# index as the rest of the table instead of creating new style xfs.
$ws.Range("C2:K2").Copy()
$ws.Range("C3:K4").PasteSpecial(-4122)

# Row 3: L1_DoubleMu0
$ws.Range("A3").Value = 100
$ws.Range("B3").Value = "L1_DoubleMu0"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 101
$ws.Range("E3").Value = 102
$ws.Range("F3").Value = 103
$ws.Range("G3").Value = 104
$ws.Range("H3").Value = 105
$ws.Range("I3").Value = 106
$ws.Range("J3").Value = 107
$ws.Range("K3").Value = 108

# Row 4: L1_DoubleMu0_SQ
$ws.Range("A4").Value = 110
$ws.Range("B4").Value = "L1_DoubleMu0_SQ"
$ws.Range("C4").Value = 110
$ws.Range("D4").Value = 111
$ws.Range("E4").Value = 112
$ws.Range("F4").Value = 113
$ws.Range("G4").Value = 114
$ws.Range("H4").Value = 115
$ws.Range("I4").Value = 116
$ws.Range("J4").Value = 117
$ws.Range("K4").Value = 118

# Row 5: L1_SingleMu28 (previously row 4, now moved above L1_SingleMuOpen).
# Index value changes 500 -> 40; the C:K measurements are unchanged.
$ws.Range("A5").Value = 40
$ws.Range("B5").Value = "L1_SingleMu28"
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 100
$ws.Range("E5").Value = 90
$ws.Range("F5").Value = 80
$ws.Range("G5").Value = 70
$ws.Range("H5").Value = 60
$ws.Range("I5").Value = 50
$ws.Range("J5").Value = 40
$ws.Range("K5").Value = 30

# Row 6: L1_SingleMuOpen (previously row 3, now moved below L1_SingleMu28).
# Values (including Index = 4) are unchanged from before.
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "L1_SingleMuOpen"
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 30
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = 6
$ws.Range("J6").Value = 7
$ws.Range("K6").Value = 8

# Row heights for the new / shifted data rows
$ws.Rows.Item(3).RowHeight = 14.9
$ws.Rows.Item(4).RowHeight = 14.9
$ws.Rows.Item(5).RowHeight = 14.2
$ws.Rows.Item(6).RowHeight = 14.2

# The sheet's used range now extends all the way down to the very last
# two rows of the worksheet (this mirrors the bigger `dimension` seen in
# the edited workbook). Touch them minimally so the used range grows,
# and give them the same row height found in the target file.
$lastRow = $ws.Rows.Count
$ws.Cells.Item($lastRow, 1).NumberFormat = "General"
$ws.Cells.Item($lastRow - 1, 1).NumberFormat = "General"
$ws.Rows.Item($lastRow).RowHeight = 12.8
$ws.Rows.Item($lastRow - 1).RowHeight = 12.8

# ------------------------------------------------------------------
# Conditional formatting: extend the existing 3 rules (currently
# applied to C2:K3) so that they cover the whole, now-bigger table
# body (C2:K6), then get rid of the redundant per-cell duplicates
# that used to cover C4:K4.
# ------------------------------------------------------------------
$headerCf = $ws.Range("C2:K3")
for ($i = 1; $i -le $headerCf.FormatConditions.Count; $i++) {
    $fc = $headerCf.FormatConditions.Item($i)
    $fc.ModifyAppliesToRange($ws.Range("C2:K6"))
}
$ws.Range("C4:K4").FormatConditions.Delete()

# Move the active selection the same way it ended up in the edited file
$ws.Range("A3").Select()
